$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the hyperlink that was on F22, then clear out the phone/email cells
# that were added (E22, F22) for "Đỗ Hoàng Băng Tâm".
$ws.Range("F22").Hyperlinks.Delete()
$ws.Range("E22:F22").Clear()

# Restore the original roster name in D22.
$ws.Range("D22").Value = "Nguyễn Duy Cường"

# Drop the now-unused built-in "Hyperlink" cell style.
$wb.Styles.Item("Hyperlink").Delete()

# Restore column F's width (21.333333333333332 in the ColumnWidth unit
# lands on the same rounded internal bucket as the original file's
# 22.140625 character-width value).
$ws.Columns.Item(6).ColumnWidth = 21.333333333333332

# Restore the previously selected cell.
$ws.Range("E23").Select()
